$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the account 005338054 / ELAINE / 100075.56 row entirely ---
$elaine = $ws.Columns(1).Find("005338054")
if ($elaine -ne $null) {
    $ws.Rows($elaine.Row).Delete()
}

# --- 2) Insert a new row for account 004350197 / GISELA / 13513.45
#        immediately before the existing 004222784 / RAFAEL row ---
$rafael = $ws.Columns(1).Find("004222784")
$newRow = $rafael.Row
$ws.Rows($newRow).Insert()

# Column A holds zero-padded account numbers stored as text; force text so the
# leading zeros are preserved, then drop the helper number format so the cell
# is left with the default (unstyled) look used by the rest of the sheet.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "004350197"
$ws.Range("A" + $newRow).ClearFormats()
$ws.Range("B" + $newRow).Value = "GISELA"
$ws.Range("C" + $newRow).Value = 13513.45

# --- 3) Remove the account 004479965 / DIEGO / 5000 row entirely ---
$diego = $ws.Columns(1).Find("004479965")
if ($diego -ne $null) {
    $ws.Rows($diego.Row).Delete()
}
